$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Core")

$ws.Range("A2").Value = $ws.Range("B2").Value2
$ws.Range("A3").Value = $ws.Range("B3").Value2
$ws.Range("A4").Value = $ws.Range("B4").Value2
$ws.Range("A5").Value = $ws.Range("B5").Value2

# Row 5's visibilityStatus cell (R5) used a slightly different style (s="2") than
# the same column in rows 2-4 (s="1"). Normalize it by copying formats from R2.
$ws.Range("R2").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
